# Natmi following Dr Hou advice
#
# Updates the Gpc3 -> Igf1r LR-pairs sheet: the Gpc3/Igf1r ligand-receptor
# edge stats are recomputed to include a new "M2" sending/target cluster
# alongside the existing ECs / FAPs / sCs clusters, giving every
# Sending-cluster x Target-cluster combination across all four clusters
# (16 data rows total, rows 2-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Gpc3", "Igf1r", "ECs", 3, 1, 3.907455, 11.722365, 0.02724526174224166, 0.02724526174224165, 3, 1, 23.18520366666667, 69.555611, 0.4216200689608106, 0.4216200689608105, 90.595139993335, 815.356259940015, 0.01148714913461926, 0.01148714913461926),
  @("ECs", "Gpc3", "Igf1r", "FAPs", 3, 1, 3.907455, 11.722365, 0.02724526174224166, 0.02724526174224165, 3, 1, 11.56543033333333, 34.696291, 0.210315924104302, 0.2103159241043019, 45.191398583135, 406.722587248215, 0.005730112400783138, 0.005730112400783137),
  @("ECs", "Gpc3", "Igf1r", "M2", 3, 1, 3.907455, 11.722365, 0.02724526174224166, 0.02724526174224165, 3, 1, 2.096289333333333, 6.288868, 0.03812076296541245, 0.03812076296541244, 8.191156236980001, 73.72040613282, 0.001038610164806614, 0.001038610164806614),
  @("ECs", "Gpc3", "Igf1r", "sCs", 3, 1, 3.907455, 11.722365, 0.02724526174224166, 0.02724526174224165, 3, 1, 18.14382633333333, 54.431479, 0.3299432439694752, 0.3299432439694752, 70.896184925315, 638.0656643278351, 0.008989390042032646, 0.008989390042032646),
  @("FAPs", "Gpc3", "Igf1r", "ECs", 3, 1, 87.66866533333332, 263.005996, 0.6112816996228118, 0.6112816996228116, 3, 1, 23.18520366666667, 69.555611, 0.4216200689608106, 0.4216200689608105, 2032.615860938173, 18293.54274844355, 0.2577286323494514, 0.2577286323494513),
  @("FAPs", "Gpc3", "Igf1r", "FAPs", 3, 1, 87.66866533333332, 263.005996, 0.6112816996228118, 0.6112816996228116, 3, 1, 11.56543033333333, 34.696291, 0.210315924104302, 0.2103159241043019, 1013.925841328982, 9125.332571960837, 0.12856227554422, 0.12856227554422),
  @("FAPs", "Gpc3", "Igf1r", "M2", 3, 1, 87.66866533333332, 263.005996, 0.6112816996228118, 0.6112816996228116, 3, 1, 2.096289333333333, 6.288868, 0.03812076296541245, 0.03812076296541244, 183.7788880058364, 1654.009992052528, 0.02330252477641566, 0.02330252477641565),
  @("FAPs", "Gpc3", "Igf1r", "sCs", 3, 1, 87.66866533333332, 263.005996, 0.6112816996228118, 0.6112816996228116, 3, 1, 18.14382633333333, 54.431479, 0.3299432439694752, 0.3299432439694752, 1590.64503868312, 14315.80534814808, 0.2016882669527248, 0.2016882669527248),
  @("M2", "Gpc3", "Igf1r", "ECs", 2, 0.6666666666666666, 0.06961200000000001, 0.208836, 0.0004853791433045106, 0.0004853791433045105, 3, 1, 23.18520366666667, 69.555611, 0.4216200689608106, 0.4216200689608105, 1.613968397644, 14.525715578796, 0.0002046455878721869, 0.0002046455878721868),
  @("M2", "Gpc3", "Igf1r", "FAPs", 2, 0.6666666666666666, 0.06961200000000001, 0.208836, 0.0004853791433045106, 0.0004853791433045105, 3, 1, 11.56543033333333, 34.696291, 0.210315924104302, 0.2103159241043019, 0.8050927363640001, 7.245834627276001, 0.0001020829630650426, 0.0001020829630650425),
  @("M2", "Gpc3", "Igf1r", "M2", 2, 0.6666666666666666, 0.06961200000000001, 0.208836, 0.0004853791433045106, 0.0004853791433045105, 3, 1, 2.096289333333333, 6.288868, 0.03812076296541245, 0.03812076296541244, 0.145926893072, 1.313342037648, 0.00001850302327026621, 0.0000185030232702662),
  @("M2", "Gpc3", "Igf1r", "sCs", 2, 0.6666666666666666, 0.06961200000000001, 0.208836, 0.0004853791433045106, 0.0004853791433045105, 3, 1, 18.14382633333333, 54.431479, 0.3299432439694752, 0.3299432439694752, 1.263028038716, 11.367252348444, 0.000160147569097015, 0.000160147569097015),
  @("sCs", "Gpc3", "Igf1r", "ECs", 3, 1, 51.77204933333334, 155.316148, 0.3609876594916421, 0.360987659491642, 3, 1, 23.18520366666667, 69.555611, 0.4216200689608106, 0.4216200689608105, 1200.345508034047, 10803.10957230643, 0.1521996418888678, 0.1521996418888677),
  @("sCs", "Gpc3", "Igf1r", "FAPs", 3, 1, 51.77204933333334, 155.316148, 0.3609876594916421, 0.360987659491642, 3, 1, 11.56543033333333, 34.696291, 0.210315924104302, 0.2103159241043019, 598.7660297785632, 5388.894268007069, 0.07592145319623381, 0.07592145319623378),
  @("sCs", "Gpc3", "Igf1r", "M2", 3, 1, 51.77204933333334, 155.316148, 0.3609876594916421, 0.360987659491642, 3, 1, 2.096289333333333, 6.288868, 0.03812076296541245, 0.03812076296541244, 108.5291947822738, 976.762753040464, 0.01376112500091991, 0.0137611250009199),
  @("sCs", "Gpc3", "Igf1r", "sCs", 3, 1, 51.77204933333334, 155.316148, 0.3609876594916421, 0.360987659491642, 3, 1, 18.14382633333333, 54.431479, 0.3299432439694752, 0.3299432439694752, 939.3430720247658, 8454.087648222892, 0.1191054394056207, 0.1191054394056207),
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $row[$j]
  }
}
